# Atualização dos gráficos 26082020
$wb = $excel.ActiveWorkbook

# --- Sheet "Mensal": update the last row (row 14) with the new August data point ---
$wsMensal = $wb.Worksheets.Item("Mensal")
$wsMensal.Cells.Item(14, 1).Value = 44066
$wsMensal.Cells.Item(14, 2).Value = 268.68
$wsMensal.Cells.Item(14, 4).Value = 7.96

# --- Sheet "Diario": append the new daily rows (384-390) ---
$wsDiario = $wb.Worksheets.Item("Diario")

$newRows = @(
    @(44060, 304.99, 248.87, 22.55),
    @(44061, 293.67, 248.87, 18),
    @(44062, 361.01, 248.87, 45.06),
    @(44063, 332.26, 248.87, 33.51),
    @(44064, 299.18, 248.87, 20.22),
    @(44065, 282.36, 248.87, 13.46),
    @(44066, 272.22, 248.87, 9.380000000000001)
)

$startRow = 384
$lastExistingRow = $startRow - 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Copy the date-column format (bold, bordered, centered, yyyy-mm-dd)
    # from the last pre-existing row so the new rows match the sheet's style.
    $wsDiario.Cells.Item($lastExistingRow, 1).Copy()
    $wsDiario.Cells.Item($r, 1).PasteSpecial(-4122)

    $wsDiario.Cells.Item($r, 1).Value = $data[0]
    $wsDiario.Cells.Item($r, 2).Value = $data[1]
    $wsDiario.Cells.Item($r, 3).Value = $data[2]
    $wsDiario.Cells.Item($r, 4).Value = $data[3]
}
